$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need NumberFormat forced to
# Text (@) before assignment, otherwise Excel auto-converts the string into a
# numeric cell (losing formatting like trailing zeros / multi-dot grouping).
# The cell Style is reset back to Normal afterwards so the cell keeps its
# original (unstyled) look - only the underlying stored value/type changes.
function Set-TextValue($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "26.471.34"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.700.92"
$ws.Range("E3").Value = "  +1.00%  "
Set-TextValue $ws "D5" "219.87"
$ws.Range("E5").Value = "  +0.84%  "
Set-TextValue $ws "D6" "0.5483"
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("E7").Value = "  +0.29%  "
Set-TextValue $ws "D8" "0.2753"
$ws.Range("E8").Value = "  +1.86%  "
Set-TextValue $ws "D9" "0.06463"
$ws.Range("E9").Value = "  +0.82%  "
Set-TextValue $ws "D10" "22.09"
$ws.Range("E10").Value = "  +0.52%  "
Set-TextValue $ws "D11" "0.07689"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "1.704.56"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E14").Value = "  +0.71%  "
Set-TextValue $ws "D15" "0.000008411"
$ws.Range("E15").Value = "  -0.44%  "
Set-TextValue $ws "D16" "65.93"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "26.523.92"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("E18").Value = "  +0.71%  "
Set-TextValue $ws "D19" "1.010"
$ws.Range("E19").Value = "  +0.31%  "
Set-TextValue $ws "D20" "11.00"
Set-TextValue $ws "D21" "191.80"
$ws.Range("E21").Value = "  +1.40%  "
Set-TextValue $ws "D22" "6.268"
$ws.Range("E22").Value = "  +1.17%  "
Set-TextValue $ws "D23" "1.011"
$ws.Range("E23").Value = "  +0.33%  "
Set-TextValue $ws "D24" "149.18"
$ws.Range("E24").Value = "  +3.36%  "
Set-TextValue $ws "D25" "0.1324"
$ws.Range("E25").Value = "  +7.21%  "
Set-TextValue $ws "D26" "7.922"
$ws.Range("E26").Value = "  +2.91%  "
Set-TextValue $ws "D27" "15.83"
$ws.Range("E27").Value = "  +0.33%  "
Set-TextValue $ws "D28" "0.06267"
$ws.Range("E28").Value = "  -5.70%  "
Set-TextValue $ws "D29" "1.381"
$ws.Range("E29").Value = "  +2.52%  "
Set-TextValue $ws "D31" "3.620"
$ws.Range("E31").Value = "  +1.54%  "
Set-TextValue $ws "D32" "3.607"
$ws.Range("E32").Value = "  +1.03%  "
Set-TextValue $ws "D33" "1.694"
$ws.Range("E33").Value = "  +1.88%  "
Set-TextValue $ws "D34" "1.042"
$ws.Range("E34").Value = "  +1.41%  "
Set-TextValue $ws "D35" "0.6191"
$ws.Range("E35").Value = "  -0.29%  "
Set-TextValue $ws "D36" "2.413"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("E37").Value = "  +2.54%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("D39").Value = "1.120.10"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("E40").Value = "  -3.76%  "
Set-TextValue $ws "D41" "0.8793"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.29%  "
Set-TextValue $ws "D43" "101.37"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "1.853.36"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("E45").Value = "  -0.54%  "
Set-TextValue $ws "D46" "57.65"
$ws.Range("E46").Value = "  +1.72%  "
Set-TextValue $ws "D47" "8.259"
$ws.Range("E47").Value = "  +1.39%  "
Set-TextValue $ws "D49" "0.05290"
$ws.Range("E49").Value = "  +0.37%  "
Set-TextValue $ws "D50" "6.163"
$ws.Range("E50").Value = "  +2.34%  "
Set-TextValue $ws "D51" "0.4306"
$ws.Range("E51").Value = "  +0.09%  "
